# Apply the data-content change described by the diff to the "devices" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# Clear the old value in A2 ("Chrome") and replace it with the two new
# cells used further along the row: I2 = "Avner", J2 = "MobileOS".
$ws.Range("A2").ClearContents()
$ws.Range("I2").Value = "Avner"
$ws.Range("J2").Value = "MobileOS"

# Update the active selection to match the post-edit cursor position (J3).
$ws.Activate()
$ws.Range("J3").Select()
